$wb = $excel.ActiveWorkbook
$count = $wb.Worksheets.Count
$last = $wb.Worksheets.Item($count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $last)
$newSheet.Name = "AugmentList"
$newSheet.Range("A1").Value = "Profession Augments"
$newSheet.Range("A7").Value = "Income Bonus"
$newSheet.Range("A8").Value = "Token Bonus"
$newSheet.Range("A9").Value = "Specialize Profession"
$newSheet.Range("A5").Value = "Augment Name"
$newSheet.Range("B5").Value = "Bonus Amount Minimum"
$newSheet.Range("A6").Value = "Expierence Bonus"
